$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-11 (Q0..Q9), columns A-G
$data = @(
    @("Q0", 0.08780261824621656, 0.7725920269955064, 2.162285864834432, 1.470471307042212, 1.484067398835772, 46),
    @("Q1", 0.2407776619532973, 0.7527870312651909, 2.110441622564052, 1.452735909435728, 1.44883219668836, 45),
    @("Q2", 0.1187069702962777, 0.783804746145044, 2.289524756223591, 1.513117561930861, 1.525893334249272, 44),
    @("Q3", 0.2312605598051813, 0.8154870963700255, 2.323257756554268, 1.524223656998627, 1.524407623678418, 43),
    @("Q4", 0.1543123306983461, 0.8126759977456647, 2.40121625977409, 1.549585834916572, 1.5605734288761, 42),
    @("Q5", 0.2357582618036748, 0.8352230691908547, 2.420446362568442, 1.555778378358705, 1.556915545132023, 41),
    @("Q6", 0.1764347851107931, 0.8312978488293193, 2.481252276672181, 1.575199122864211, 1.585227652752555, 40),
    @("Q7", 0.2294621647362194, 0.8569143726678996, 2.527131915378803, 1.589695541724516, 1.593611276175175, 39),
    @("Q8", 0.1874282125274516, 0.8425231446840433, 2.595659384232942, 1.611105019616332, 1.621645326275095, 38),
    @("Q9", 0.1764013765962084, 0.8338957443602626, 2.565872115365881, 1.60183398495783, 1.614052205362943, 37)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
}

# Apply the same formatting (bold, centered, bordered) as the existing
# A2:A5 cells to the new A6:A11 cells
$ws.Range("A2").Copy()
$ws.Range("A6:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
